# Dodat ispis poruke sa brojem osvojenih poena
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @("afkjl", "smv", "ksdm", "cxv", "21"),
    @("fdgdfg", "dfg", "fdg", "dfg", "11"),
    @("v", "v", "xv", "", "21"),
    @("f", "xcx", "xb", "", "15"),
    @("sdg", "xb", "vcb", "", "18")
)

$startRow = 61
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 5; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
